$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Sheets.Item("Overview")
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsDeDe = $wb.Sheets.Item("de-de")

# Latest HO Xliff Generate Date (Overview!G) / Latest Handoff Datetime (de-de!H)
# both originally shared the same text and are updated together.
$wsOverview.Range("G9:G14").Value = "2016-09-04 06:23:58"
$wsDeDe.Range("H9:H14").Value = "2016-09-04 06:23:58"

# Latest Handoff Datetime (zh-cn!H)
$wsZhCn.Range("H9:H14").Value = "2016-09-04 06:23:53"

# Priority column (E) set to "ht" for rows 9-14 on both locale sheets
$wsZhCn.Range("E9:E14").Value = "ht"
$wsDeDe.Range("E9:E14").Value = "ht"
